# fdo#75168 regression test fixture:
# add a second worksheet ("Sheet2") after "Sheet1" containing a small
# 3-row/2-column table plus two "expression" type conditional formatting
# rules (one per column), then leave Sheet2 as the active sheet with
# C1 selected.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheet right after Sheet1 and name it.
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

# Data for the new sheet.
$ws2.Range("A1").Value = 2
$ws2.Range("B1").Value = 2
$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = 1
$ws2.Range("A3").Value = 4
$ws2.Range("B3").Value = 3

# Column A: highlight cells where A1<>1 (expression rule).
$fcA = $ws2.Range("A1:A3").FormatConditions.Add(2, 0, "A1<>1")

# Column B: highlight cells where B1=1 (expression rule), added after
# column A's rule so it ends up with the higher (more recent) priority.
$fcB = $ws2.Range("B1:B3").FormatConditions.Add(2, 0, "B1=1")

# Match the priority ordering of the reference workbook: the B rule
# (added last) is priority 1 (highest), the A rule is priority 2.
$fcA.Priority = 2
$fcB.Priority = 1

# Leave the selection on C1 and Sheet2 as the active sheet/tab.
$ws2.Range("C1").Select()
